$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of "Price" column values in this diff are plain numeric-looking
# strings (e.g. "582.08"). The workbook stores every cell in columns D/E as
# literal text, so force just those specific cells to Text format first -
# otherwise Excel auto-converts the assigned string into a Number, which can
# silently change the displayed text (e.g. "1.00" -> 1).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "67.841.00"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "3.248.19"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D5").Value = "582.08"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").Value = "183.03"
$ws.Range("E6").Value = "  +4.14%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.597"
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("D9").Value = "0.135"
$ws.Range("E9").Value = "  +4.57%  "
$ws.Range("D10").Value = "6.69"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").Value = "0.416"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").Value = "3.809.97"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "28.69"
$ws.Range("E14").Value = "  +3.38%  "
$ws.Range("D15").Value = "67.848.72"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  +2.36%  "
$ws.Range("D17").Value = "3.251.05"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").Value = "13.56"
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("D20").Value = "380.08"
$ws.Range("E20").Value = "  +3.34%  "
$ws.Range("D21").Value = "7.64"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "71.39"
$ws.Range("E23").Value = "  +1.86%  "
$ws.Range("D24").Value = "0.513"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("D26").Value = "9.93"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("E27").Value = "  +2.26%  "
$ws.Range("D29").Value = "1.98"
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("D30").Value = "5.67"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("D32").Value = "7.07"
$ws.Range("E32").Value = "  +4.80%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +2.95%  "
$ws.Range("E35").Value = "  +4.28%  "
$ws.Range("D36").Value = "162.27"
$ws.Range("E36").Value = "  -6.70%  "
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").Value = "26.46"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("E40").Value = "  +5.13%  "
$ws.Range("E41").Value = "  +7.17%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "25.44"
$ws.Range("E43").Value = "  +3.87%  "
$ws.Range("D44").Value = "41.18"
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("D45").Value = "345.98"
$ws.Range("E45").Value = "  +4.24%  "
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("D47").Value = "2.621.08"
$ws.Range("E47").Value = "  -3.09%  "
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").Value = "0.103"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").Value = "0.992"
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("E51").Value = "  +2.72%  "
